$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 1641.6666
$ws.Cells.Item(28, 9).Value = 2360
$ws.Cells.Item(28, 10).Value = 1128.5714
$ws.Cells.Item(28, 11).Value = 2360
$ws.Cells.Item(28, 12).Value = 1128.5714
$ws.Cells.Item(28, 13).Value = -1875
$ws.Cells.Item(28, 14).Value = -2098.5714

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 1630.5238
$ws.Cells.Item(112, 10).Value = 1687.05
$ws.Cells.Item(112, 12).Value = 5061.15
$ws.Cells.Item(112, 14).Value = -7277.15

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3336293.5
$ws.Cells.Item(137, 9).Value = 5557894.5
$ws.Cells.Item(137, 10).Value = 3891.6667
$ws.Cells.Item(137, 11).Value = 16673683.5
$ws.Cells.Item(137, 12).Value = 11675.0001
$ws.Cells.Item(137, 13).Value = -16671133.5
$ws.Cells.Item(137, 14).Value = -16775.0001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1764.52
$ws.Cells.Item(138, 10).Value = 2016.6747
$ws.Cells.Item(138, 12).Value = 6050.024100000001
$ws.Cells.Item(138, 14).Value = -16330.0241

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1507.4286
$ws.Cells.Item(2, 9).Value = 1302.6875
$ws.Cells.Item(2, 11).Value = 1302.6875
$ws.Cells.Item(2, 13).Value = -1189.6875

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4185.2583
$ws.Cells.Item(45, 9).Value = 4156.9614
$ws.Cells.Item(45, 11).Value = 4156.9614
$ws.Cells.Item(45, 13).Value = -3779.9614

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3226
$ws.Cells.Item(63, 10).Value = 4999.5
$ws.Cells.Item(63, 12).Value = 4999.5
$ws.Cells.Item(63, 14).Value = -6371.5

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3226
$ws.Cells.Item(66, 10).Value = 4999.5
$ws.Cells.Item(66, 12).Value = 24997.5
$ws.Cells.Item(66, 14).Value = -31861.5

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 1507.4286
$ws.Cells.Item(116, 9).Value = 1302.6875
$ws.Cells.Item(116, 11).Value = 1302.6875
$ws.Cells.Item(116, 13).Value = 991.3125

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 202971.2
$ws.Cells.Item(132, 9).Value = 145958.86
$ws.Cells.Item(132, 10).Value = 336000
$ws.Cells.Item(132, 11).Value = 437876.58
$ws.Cells.Item(132, 12).Value = 1008000
$ws.Cells.Item(132, 13).Value = -435346.58
$ws.Cells.Item(132, 14).Value = -1013060

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1507.4286
$ws.Cells.Item(3, 9).Value = 1302.6875
$ws.Cells.Item(3, 11).Value = 1302.6875
$ws.Cells.Item(3, 13).Value = -1188.6875

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 984.5454999999999
$ws.Cells.Item(20, 9).Value = 889
$ws.Cells.Item(20, 11).Value = 889
$ws.Cells.Item(20, 13).Value = -642

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1169.1111
$ws.Cells.Item(94, 9).Value = 1007.4286
$ws.Cells.Item(94, 10).Value = 1735
$ws.Cells.Item(94, 11).Value = 1007.4286
$ws.Cells.Item(94, 12).Value = 1735
$ws.Cells.Item(94, 13).Value = -556.4286
$ws.Cells.Item(94, 14).Value = -2637

# BSM row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(97, 8).Value = 2758.4285
$ws.Cells.Item(97, 9).Value = 2758.4285
$ws.Cells.Item(97, 11).Value = 2758.4285
$ws.Cells.Item(97, 13).Value = -1767.4285

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2751.625
$ws.Cells.Item(107, 10).Value = 2835.3333
$ws.Cells.Item(107, 12).Value = 2835.3333
$ws.Cells.Item(107, 14).Value = -6675.3333

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4240.4614
$ws.Cells.Item(134, 9).Value = 3621.8708
$ws.Cells.Item(134, 10).Value = 6637.5
$ws.Cells.Item(134, 11).Value = 10865.6124
$ws.Cells.Item(134, 12).Value = 19912.5
$ws.Cells.Item(134, 13).Value = -8330.6124
$ws.Cells.Item(134, 14).Value = -24982.5

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2682.7856
$ws.Cells.Item(31, 9).Value = 1250.5
$ws.Cells.Item(31, 11).Value = 1250.5
$ws.Cells.Item(31, 13).Value = -955.5

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2682.7856
$ws.Cells.Item(34, 9).Value = 1250.5
$ws.Cells.Item(34, 11).Value = 1250.5
$ws.Cells.Item(34, 13).Value = -1048.5

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 9).Value = 40003860
$ws.Cells.Item(58, 10).Value = 1999.25
$ws.Cells.Item(58, 11).Value = 40003860
$ws.Cells.Item(58, 12).Value = 1999.25
$ws.Cells.Item(58, 13).Value = -40003657
$ws.Cells.Item(58, 14).Value = -2405.25

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 1271.4
$ws.Cells.Item(107, 9).Value = 521.375
$ws.Cells.Item(107, 10).Value = 2128.5715
$ws.Cells.Item(107, 11).Value = 521.375
$ws.Cells.Item(107, 12).Value = 2128.5715
$ws.Cells.Item(107, 13).Value = 1398.625
$ws.Cells.Item(107, 14).Value = -5968.5715

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 9).Value = 40003860
$ws.Cells.Item(136, 10).Value = 1999.25
$ws.Cells.Item(136, 11).Value = 120011580
$ws.Cells.Item(136, 12).Value = 5997.75
$ws.Cells.Item(136, 13).Value = -120009030
$ws.Cells.Item(136, 14).Value = -11097.75

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 600
$ws.Cells.Item(92, 10).Value = 900
$ws.Cells.Item(92, 12).Value = 2700
$ws.Cells.Item(92, 14).Value = -5196

# CUL row 127
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(127, 8).Value = 677.6667
$ws.Cells.Item(127, 10).Value = 677.6667
$ws.Cells.Item(127, 12).Value = 2033.0001
$ws.Cells.Item(127, 14).Value = -11953.0001

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 956.01514
$ws.Cells.Item(131, 9).Value = 800
$ws.Cells.Item(131, 10).Value = 960.8905999999999
$ws.Cells.Item(131, 11).Value = 2400
$ws.Cells.Item(131, 12).Value = 2882.6718
$ws.Cells.Item(131, 13).Value = 2640
$ws.Cells.Item(131, 14).Value = -12962.6718

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2156.25
$ws.Cells.Item(132, 9).Value = 2156.25
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 19406.25
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -16876.25
$ws.Cells.Item(132, 14).ClearContents()

# CUL row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(133, 8).Value = 4341
$ws.Cells.Item(133, 9).Value = 3521.6667
$ws.Cells.Item(133, 10).Value = 6799
$ws.Cells.Item(133, 11).Value = 10565.0001
$ws.Cells.Item(133, 12).Value = 20397
$ws.Cells.Item(133, 13).Value = -5505.000100000001
$ws.Cells.Item(133, 14).Value = -30517

# CUL row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 5454.4443
$ws.Cells.Item(141, 9).Value = 3181.6667
$ws.Cells.Item(141, 11).Value = 9545.000100000001
$ws.Cells.Item(141, 13).Value = -4365.000100000001

# GSM row 43
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 7758.5
$ws.Cells.Item(43, 9).Value = 517
$ws.Cells.Item(43, 10).Value = 15000
$ws.Cells.Item(43, 11).Value = 517
$ws.Cells.Item(43, 12).Value = 15000
$ws.Cells.Item(43, 13).Value = -366
$ws.Cells.Item(43, 14).Value = -15302

# GSM row 46
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 14).ClearContents()

# GSM row 57
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 17800
$ws.Cells.Item(57, 9).Value = 17800
$ws.Cells.Item(57, 11).Value = 17800
$ws.Cells.Item(57, 13).Value = -16980

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 20840760
$ws.Cells.Item(80, 9).Value = 10077.214
$ws.Cells.Item(80, 10).Value = 50003716
$ws.Cells.Item(80, 11).Value = 10077.214
$ws.Cells.Item(80, 12).Value = 50003716
$ws.Cells.Item(80, 13).Value = -9079.214
$ws.Cells.Item(80, 14).Value = -50005712

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 20840760
$ws.Cells.Item(83, 9).Value = 10077.214
$ws.Cells.Item(83, 10).Value = 50003716
$ws.Cells.Item(83, 11).Value = 50386.07
$ws.Cells.Item(83, 12).Value = 250018580
$ws.Cells.Item(83, 13).Value = -45394.07
$ws.Cells.Item(83, 14).Value = -250028564

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 1833.1818
$ws.Cells.Item(107, 9).Value = 1551.8
$ws.Cells.Item(107, 10).Value = 2067.6667
$ws.Cells.Item(107, 11).Value = 1551.8
$ws.Cells.Item(107, 12).Value = 2067.6667
$ws.Cells.Item(107, 13).Value = 368.2
$ws.Cells.Item(107, 14).Value = -5907.6667

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1976
$ws.Cells.Item(113, 9).Value = 1750.25
$ws.Cells.Item(113, 10).Value = 2105
$ws.Cells.Item(113, 11).Value = 1750.25
$ws.Cells.Item(113, 12).Value = 2105
$ws.Cells.Item(113, 13).Value = 419.75
$ws.Cells.Item(113, 14).Value = -6445

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 80517.38
$ws.Cells.Item(132, 9).Value = 70221.2
$ws.Cells.Item(132, 10).Value = 94557.63
$ws.Cells.Item(132, 11).Value = 210663.6
$ws.Cells.Item(132, 12).Value = 283672.89
$ws.Cells.Item(132, 13).Value = -208133.6
$ws.Cells.Item(132, 14).Value = -288732.89

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 20456.5
$ws.Cells.Item(136, 10).Value = 20456.5
$ws.Cells.Item(136, 12).Value = 61369.5
$ws.Cells.Item(136, 14).Value = -66469.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 4142.7144
$ws.Cells.Item(61, 9).Value = 7000
$ws.Cells.Item(61, 10).Value = 2999.8
$ws.Cells.Item(61, 11).Value = 7000
$ws.Cells.Item(61, 12).Value = 2999.8
$ws.Cells.Item(61, 13).Value = -6798
$ws.Cells.Item(61, 14).Value = -3403.8

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 4142.7144
$ws.Cells.Item(113, 9).Value = 7000
$ws.Cells.Item(113, 10).Value = 2999.8
$ws.Cells.Item(113, 11).Value = 7000
$ws.Cells.Item(113, 12).Value = 2999.8
$ws.Cells.Item(113, 13).Value = -4830
$ws.Cells.Item(113, 14).Value = -7339.8

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3296.5293
$ws.Cells.Item(122, 9).Value = 3011.75
$ws.Cells.Item(122, 10).Value = 3980
$ws.Cells.Item(122, 11).Value = 9035.25
$ws.Cells.Item(122, 12).Value = 11940
$ws.Cells.Item(122, 13).Value = -6585.25
$ws.Cells.Item(122, 14).Value = -16840

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 884.1429000000001
$ws.Cells.Item(113, 9).Value = 568.875
$ws.Cells.Item(113, 10).Value = 1078.1538
$ws.Cells.Item(113, 11).Value = 1706.625
$ws.Cells.Item(113, 12).Value = 3234.4614
$ws.Cells.Item(113, 13).Value = 463.375
$ws.Cells.Item(113, 14).Value = -7574.4614

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 85960.75
$ws.Cells.Item(136, 9).Value = 58580.11
$ws.Cells.Item(136, 10).Value = 168102.67
$ws.Cells.Item(136, 11).Value = 175740.33
$ws.Cells.Item(136, 12).Value = 504308.01
$ws.Cells.Item(136, 13).Value = -173190.33
$ws.Cells.Item(136, 14).Value = -509408.01

# WVR row 137
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(137, 8).Value = 49800
$ws.Cells.Item(137, 10).Value = 49800
$ws.Cells.Item(137, 12).Value = 49800
$ws.Cells.Item(137, 14).Value = -60000
